$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="27.209.80"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("E2").Value = '  +1.22%  '

$ws.Range("D3").Formula = '="1.643.29"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("E3").Value = '  +0.05%  '

$ws.Range("D5").Formula = '="217.22"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)

$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("E8").Value = '  +0.95%  '

$ws.Range("E9").Value = '  +0.79%  '

$ws.Range("E10").Value = '  +1.36%  '

$ws.Range("D11").Formula = '="0.0848"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)

$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("D12").Formula = '="1.873.10"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)

$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").Formula = '="1.646.06"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)

$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("E15").Value = '  +2.86%  '

$ws.Range("D16").Formula = '="67.42"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)

$ws.Range("E16").Value = '  +1.68%  '

$ws.Range("D17").Formula = '="27.190.56"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)

$ws.Range("D18").Formula = '="0.0"&UNICHAR(8323)&"0740"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)

$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("D19").Formula = '="219.02"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)

$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").Formula = '="6.88"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)

$ws.Range("E21").Value = '  +3.78%  '

$ws.Range("D22").Formula = '="2.57"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)

$ws.Range("E22").Value = '  +6.51%  '

$ws.Range("E23").Value = '  +0.55%  '

$ws.Range("D24").Formula = '="9.23"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)

$ws.Range("E24").Value = '  +0.44%  '

$ws.Range("D25").Formula = '="147.83"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)

$ws.Range("E25").Value = '  +1.20%  '

$ws.Range("D26").Formula = '="7.54"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)

$ws.Range("E26").Value = '  +1.83%  '

$ws.Range("E27").Value = '  -0.04%  '

$ws.Range("E28").Value = '  -0.45%  '

$ws.Range("D29").Formula = '="15.80"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)

$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").Formula = '="0.0509"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)

$ws.Range("E30").Value = '  +0.67%  '

$ws.Range("D31").Formula = '="1.18"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)

$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("E32").Value = '  +1.00%  '

$ws.Range("E33").Value = '  +1.00%  '

$ws.Range("E34").Value = '  +1.24%  '

$ws.Range("D35").Formula = '="1.275.32"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)

$ws.Range("E35").Value = '  +2.35%  '

$ws.Range("E36").Value = '  +0.95%  '

$ws.Range("E37").Value = '  +1.71%  '

$ws.Range("E38").Value = '  +1.47%  '

$ws.Range("E39").Value = '  +2.67%  '

$ws.Range("E40").Value = '  -0.09%  '

$ws.Range("E41").Value = '  +0.48%  '

$ws.Range("D42").Formula = '="2.24"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)

$ws.Range("E42").Value = '  +7.14%  '

$ws.Range("E43").Value = '  -1.00%  '

$ws.Range("D44").Formula = '="1.782.86"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").Formula = '="61.84"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)

$ws.Range("E45").Value = '  +1.75%  '

$ws.Range("E47").Value = '  +1.53%  '

$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("D50").Formula = '="7.69"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)

$ws.Range("E50").Value = '  +0.93%  '

$ws.Range("D51").Formula = '="0.0974"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$ws.Range("E51").Value = '  +0.05%  '

